$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as numbers by Excel;
# force Text number format first so the OOXML stores them as strings, matching the source data.
$numericLookingCells = @("D5","D6","D9","D10","D11","D12","D13","D14","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D39","D40","D41","D44","D45","D46","D48","D49","D50","D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = '573.73'
$ws.Range("D6").Value = '148.51'
$ws.Range("D9").Value = '0.524'
$ws.Range("D10").Value = '0.158'
$ws.Range("D11").Value = '6.06'
$ws.Range("D12").Value = '0.495'
$ws.Range("D13").Value = '0.0000264'
$ws.Range("D14").Value = '36.88'
$ws.Range("D18").Value = '7.06'
$ws.Range("D20").Value = '500.48'
$ws.Range("D21").Value = '14.72'
$ws.Range("D22").Value = '0.713'
$ws.Range("D23").Value = '15.13'
$ws.Range("D24").Value = '7.69'
$ws.Range("D25").Value = '83.97'
$ws.Range("D26").Value = '0.999'
$ws.Range("D27").Value = '8.86'
$ws.Range("D28").Value = '2.89'
$ws.Range("D29").Value = '2.15'
$ws.Range("D30").Value = '2.77'
$ws.Range("D31").Value = '27.40'
$ws.Range("D32").Value = '1.00'
$ws.Range("D34").Value = '6.15'
$ws.Range("D35").Value = '6.45'
$ws.Range("D36").Value = '54.27'
$ws.Range("D37").Value = '0.0896'
$ws.Range("D38").Value = '465.07'
$ws.Range("D39").Value = '0.0416'
$ws.Range("D40").Value = '2.97'
$ws.Range("D41").Value = '8.61'
$ws.Range("D44").Value = '2.42'
$ws.Range("D45").Value = '0.281'
$ws.Range("D46").Value = '28.25'
$ws.Range("D48").Value = '0.999'
$ws.Range("D49").Value = '0.114'
$ws.Range("D50").Value = '2.23'
$ws.Range("D51").Value = '118.61'

# Remaining text cells (already safe from numeric auto-conversion)
$ws.Range("D2").Value = '64.521.92'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '3.137.67'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("E6").Value = '  +0.84%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").Value = '3.137.68'
$ws.Range("E8").Value = '  +1.21%  '
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("E11").Value = '  -1.28%  '
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("E13").Value = '  +10.92%  '
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").Value = '3.662.21'
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("D16").Value = '64.815.14'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '3.147.24'
$ws.Range("E17").Value = '  +1.46%  '
$ws.Range("E18").Value = '  +0.66%  '
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("E23").Value = '  -0.43%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  +4.96%  '
$ws.Range("E28").Value = '  +0.83%  '
$ws.Range("E29").Value = '  +1.39%  '
$ws.Range("E30").Value = '  +3.34%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("E34").Value = '  +2.85%  '
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("E36").Value = '  -2.01%  '
$ws.Range("E37").Value = '  +7.18%  '
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("E39").Value = '  -0.82%  '
$ws.Range("E40").Value = '  +4.44%  '
$ws.Range("E41").Value = '  +1.79%  '
$ws.Range("D42").Value = '3.029.53'
$ws.Range("E42").Value = '  -1.50%  '
$ws.Range("E43").Value = '  -2.04%  '
$ws.Range("E44").Value = '  +3.99%  '
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("E46").Value = '  -1.13%  '
$ws.Range("D47").Value = '0.0₃0581'
$ws.Range("E47").Value = '  +7.93%  '
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("E49").Value = '  -0.72%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("E50").Value = '  +1.75%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E51").Value = '  +0.66%  '

Write-Host "Applied 99 cell updates"
